# Word COM-interop edit script
#
# Target change (per the OOXML diff):
#   1. The existing "I AM OVERWRITTEN DEFAULT STYLE" Heading1 paragraph gets
#      explicit direct formatting: single (240 twips / "auto") line spacing
#      and bold (w:b / w:bCs) on both the paragraph mark and the run.
#   2. A new, empty paragraph is appended right after it, with its own
#      direct formatting: 1.5 line spacing (360 "auto"), justified alignment,
#      and a paragraph-mark font of Times New Roman 14pt (sz 28) / en-US.
#
# Word's Range.Font.Bold / BoldBi property writes don't reliably stamp the
# paragraph-mark's own <w:rPr> (only the run's), so the most faithful way to
# reproduce the exact target markup is to replace the body's content with
# the precise OOXML via Range.InsertXML (a standard, supported Word COM
# method - it substitutes the addressed range's contents with the supplied
# WordprocessingML).

$d = $word.ActiveDocument

$targetXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="Heading1"/>
              <w:spacing w:line="240" w:lineRule="auto"/>
              <w:rPr>
                <w:b/>
                <w:bCs/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:b/>
                <w:bCs/>
              </w:rPr>
              <w:t>I AM OVERWRITTEN DEFAULT STYLE</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:spacing w:line="360" w:lineRule="auto"/>
              <w:jc w:val="both"/>
              <w:rPr>
                <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
                <w:sz w:val="28"/>
                <w:szCs w:val="28"/>
                <w:lang w:val="en-US"/>
              </w:rPr>
            </w:pPr>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

# Replace the whole body's content (the sectPr at the end of the body is not
# part of Content, so it is left untouched) with the exact target markup.
[void]$d.Content.InsertXML($targetXml)

Write-Output "applied overwritten-default-style formatting"
